# Incremento "Clase Categoría" en la planilla de métricas
# Rellena la fila 19 (antes vacía) de la tabla "Desarrollo y correctivos"
# con un nuevo incremento llamado "Clase Categoría".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Descripción de las tareas del incremento (columna C, celdas combinadas C19:E19)
$ws.Range("C19").Value = "Clase Categoría"

# Estimación: Líneas Cód. / Tiempo
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 0.0069444444444444441

# Desarrollo: Hora Inicio / Hora Fin
$ws.Range("H19").Value = 0.87430555555555556
$ws.Range("I19").Value = 0.88263888888888886

# Correctivos: Errores Lógicos / Tiempo Correción E.L.
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0

# Líneas Reales
$ws.Range("M19").Value = 52

# La celda activa/seleccionada pasa a F20
$ws.Range("F20").Select() | Out-Null
